$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Username / Password (bold, yellow fill - reuse existing header style) ---
$ws.Range("A1").Value = "Username"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535

$ws.Range("B1").Value = "Password"
# Copy A1's formatting onto B1 so both share the same cell style (avoids creating
# a spurious intermediate style record).
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data row: credentials, turned into live hyperlinks ---
$ws.Range("A2").Value = "Valmiki.k@intelegain.com"
$ws.Range("B2").Value = "Test@123"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Valmiki.k@intelegain.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Test@123")

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 26
$ws.Columns("B").ColumnWidth = 9.3
